$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 1121
$wsExhibition.Range("F4").Value = 1817
$wsExhibition.Range("F6").Value = 441

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1121
$wsAll.Range("F4").Value = 1817
$wsAll.Range("F7").Value = 441
